# Final upload/last changes for thesis composition:
# Add a new "Comparison to other Studies" sheet (after Sheet1) containing two
# more Cohen's Kappa contingency-table comparisons (pulled from the paper at
# https://ieeexplore.ieee.org/document/6449272/), and update the selection on
# Sheet1 to reflect the cells that were last highlighted there.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: selection moves from I13 to A2:G6 ---------------------------
[void]$ws1.Range("A2:G6").Select()

# --- New worksheet, inserted right after Sheet1 ---------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Comparison to other Studies"

# First comparison table (rows 1-7)
$ws2.Range("C1").Value = "TRUTH"

$ws2.Range("C2").Value = "Class A"
$ws2.Range("D2").Value = "Class B"

$ws2.Range("A3").Value = "PREDICT"
$ws2.Range("B3").Value = "Class A"
$ws2.Range("C3").Value = 28
$ws2.Range("D3").Value = 4
$ws2.Range("F3").Value = "Obs. Acc:"
$ws2.Range("G3").Formula = "=(C3+D4)/SUM(C3:D4)"

$ws2.Range("B4").Value = "Class B"
$ws2.Range("C4").Value = 2
$ws2.Range("D4").Value = 7
$ws2.Range("F4").Value = "Exp. Acc:"
$ws2.Range("G4").Formula = "=(SUM(C3:C4)*SUM(C3:D3) + SUM(D3:D4)*SUM(C4:D4))/(SUM(C3:D4)^2)"

$ws2.Range("F5").Value = "Kappa:"
$ws2.Range("G5").Formula = "=(G3-G4)/(1-G4)"

$ws2.Range("B7").Value = "https://ieeexplore.ieee.org/document/6449272/"

# Second comparison table (rows 10-16)
$ws2.Range("C10").Value = "TRUTH"

$ws2.Range("C11").Value = "Class A"
$ws2.Range("D11").Value = "Class B"

$ws2.Range("A12").Value = "PREDICT"
$ws2.Range("B12").Value = "Class A"
$ws2.Range("C12").Value = 14
$ws2.Range("D12").Value = 3
$ws2.Range("F12").Value = "Obs. Acc:"
$ws2.Range("G12").Formula = "=(C12+D13)/SUM(C12:D13)"

$ws2.Range("B13").Value = "Class B"
$ws2.Range("C13").Value = 3
$ws2.Range("D13").Value = 9
$ws2.Range("F13").Value = "Exp. Acc:"
$ws2.Range("G13").Formula = "=(SUM(C12:C13)*SUM(C12:D12) + SUM(D12:D13)*SUM(C13:D13))/(SUM(C12:D13)^2)"

$ws2.Range("F14").Value = "Kappa:"
$ws2.Range("G14").Formula = "=(G12-G13)/(1-G13)"

$ws2.Range("B16").Value = "https://ieeexplore.ieee.org/document/6449272/"

# Leave the new sheet's selection on D12, matching the final saved state.
[void]$ws2.Range("D12").Select()
